# Update NB3 Body (for servo drive)
#
# Row 2 of the BOM describes the "Header Breakaway 6 2 Surface Mount"
# connector. The servo-drive variant now also populates J16/J15, and the
# quantity used on the board goes from 3 to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity (column A) for the Header Breakaway part: 3 -> 5
$ws.Range("A2").Value = 5

# Designator list (column B) for the same part: add J16, J15 ahead of the
# existing J10/J12/J11 designators.
$ws.Range("B2").Value = "J16,J15,J10,J12,J11"

# Re-apply the "Normal" cell style to B2 so the edited cell carries its own
# explicit font formatting (mirrors what Excel does when the cell's content
# is retyped).
$ws.Range("B2").Style = "Normal"

# Leave the cursor where the user ended up after making the edit.
$ws.Range("C18").Select() | Out-Null
